# Apply the "Add files via upload" edit to the CLEANUP_DETAILS sheet:
# - Insert a new column before column D ("INVALID RECORDS ALGORITHM"),
#   pushing the existing "CLEANING ORDER" column (D) to E and the trailing
#   blank column (E) to F.
# - Rename the old "VALID RELATION TO " header (C1) to "VALID RECORDS ALGORITHM".
# - Give the new column D the header "INVALID RECORDS ALGORITHM".
# - Resize/cleanup the sheet view, column widths and row heights to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D. Excel shifts D->E and E->F, and copies
#    formatting from the column to the left (C), matching the target layout.
$ws.Columns("D").Insert()

# 2. Update the header row text.
$ws.Range("C1").Value = "VALID RECORDS ALGORITHM"
$ws.Range("D1").Value = "INVALID RECORDS ALGORITHM"

# 3. Match column widths: column D should be exactly as wide as column C
#    (copy the existing width rather than assigning a literal, to avoid
#    pixel-rounding drift).
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# 4. Adjust row heights to match the new, narrower wrapped text.
$ws.Rows("1").RowHeight = 29
$ws.Rows("2").RowHeight = 87
$ws.Rows("3").RowHeight = 87

# 5. Reset the view: scroll back to the top and select D3.
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("D3").Select()
